$wb = $excel.ActiveWorkbook

# Overview sheet: row 3 corresponds to b3b82a64-5468-4710-a2a7-13a2a451b96b.md
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# zh-cn sheet: row 3 Status -> Ready for handoff, Latest Handoff Datetime updated
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "2016-03-10 00:14:39"

# de-de sheet: row 3 Status -> Ready for handoff, Latest Handoff Datetime updated
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "2016-03-10 00:14:44"
